# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text must stay as text (avoid Excel
# auto-converting e.g. "87.69" into a numeric value) get NumberFormat "@"
# applied first, exactly as typing into a Text-formatted cell would.

$ws.Range("D2").Value = "40.172.54"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").Value = "2.227.68"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.19"
$ws.Range("E5").Value = "  +1.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.69"
$ws.Range("E6").Value = "  -0.92%  "

$ws.Range("E7").Value = "  -0.50%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.469"
$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.65"
$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.97"
$ws.Range("E11").Value = "  +7.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0782"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("E13").Value = "  +2.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.45"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").Value = "2.575.31"
$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.83"
$ws.Range("E16").Value = "  -1.18%  "

$ws.Range("D17").Value = "2.232.77"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("E18").Value = "  +1.21%  "

$ws.Range("D19").Value = "40.112.39"
$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.27"
$ws.Range("E21").Value = "  -4.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.78"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.80"
$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.89"
$ws.Range("E24").Value = "  +0.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.23%  "

$ws.Range("E26").Value = "  +0.65%  "

$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.31"
$ws.Range("E28").Value = "  +3.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.58"
$ws.Range("E31").Value = "  +2.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.77"
$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.97"
$ws.Range("E34").Value = "  +0.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.05"
$ws.Range("E35").Value = "  +7.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0716"
$ws.Range("E36").Value = "  -0.30%  "

$ws.Range("E37").Value = "  -3.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.75"
$ws.Range("E39").Value = "  +2.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0994"
$ws.Range("E40").Value = "  -0.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.46"
$ws.Range("E41").Value = "  -3.63%  "

$ws.Range("D42").Value = "2.090.63"
$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.72"
$ws.Range("E43").Value = "  -2.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.55"
$ws.Range("E44").Value = "  +4.97%  "

$ws.Range("E45").Value = "  +2.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0271"
$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.95"
$ws.Range("E47").Value = "  -11.16%  "

$ws.Range("E48").Value = "  +2.28%  "

$ws.Range("D49").Value = "2.448.50"
$ws.Range("E49").Value = "  +0.76%  "

$ws.Range("E50").Value = "  +2.68%  "

$ws.Range("E51").Value = "  +3.49%  "

# Rows 29/30: Toncoin and Cosmos swap places
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.31"
$ws.Range("E29").Value = "  +0.64%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.05"
$ws.Range("E30").Value = "  -6.60%  "

